$d = $word.ActiveDocument

# The word "Bayesian " carries an open reviewer comment (id=3, "To discuss").
# The edit drops the word "Bayesian" in favor of "a" and resolves/removes
# that comment entirely (text + commentRangeStart/End + commentReference +
# the matching commentsExtended/commentsExtensible/commentsIds entries).

$comment = $null
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $candidate = $d.Comments.Item($i)
    if ($candidate.Scope.Text -eq "Bayesian ") {
        $comment = $candidate
    }
}

$scope = $comment.Scope
$start = $scope.Start
$end = $scope.End

# Deleting the comment removes the comment-range markers/reference and the
# comment definition (comments.xml / commentsExtended.xml /
# commentsExtensible.xml / commentsIds.xml) but leaves the commented text
# ("Bayesian ") in place.
$comment.Delete()

# Replace the now-uncommented text with "a ".
$range = $d.Range($start, $end)
$range.Text = "a "
